$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-DateCell($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $text

    # Apply the "Times New Roman" / 28 half-points formatting used throughout
    # the schedule table. The complex-script font (w:cs) has to be applied
    # via a formatted Find/Replace first (Font.NameBi on a table-cell range
    # isn't directly settable), then the ascii/hAnsi/sz/szCs via Font.

    $rngFind = $d.Content
    $rngFind.Find.ClearFormatting()
    $rngFind.Find.Text = $text
    $rngFind.Find.Replacement.ClearFormatting()
    $rngFind.Find.Replacement.Font.NameBi = "Times New Roman"
    $rngFind.Find.Replacement.Text = $text
    $rngFind.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2, $true)

    $cell2 = $table.Cell($row, $col)
    $cell2.Range.Font.Name = "Times New Roman"
    $cell2.Range.Font.Size = 14
    $cell2.Range.Font.SizeBi = 14
}

# Row 19 (1-based Word table row) == the "Л09 / 22.11 / 20.11" row's
# follow-up row (ЛР05), whose four schedule-date cells are currently empty.
Set-DateCell $t 19 2 "30.11"
Set-DateCell $t 19 3 "03.12"
Set-DateCell $t 19 4 "27.11"
Set-DateCell $t 19 5 "01.12"

# Row 20 (1-based Word table row) == the "Л10" row.
Set-DateCell $t 20 2 "06.12"
Set-DateCell $t 20 4 "04.12"

# The "_GoBack" bookmark (Word's last-edit marker) moves from the "22.11"
# cell to the end of the newly-typed "06.12" cell, mirroring Word's own
# behaviour of tracking the most recent edit location.
if (-not $d.Bookmarks.Exists("_GoBack")) {
    $endRange = $t.Cell(20, 2).Range
    $endRange.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $endRange)
}
